# Update column C (Förändrad) values from 46061 to 46062 for rows 2-23
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
